# Insert a new weekly price-report row for "Terminal La Palmera de La Serena"
# (Cebollín) right before the current row 246, pushing every following row
# down by one (old row 246 -> 247, ..., old row 281 -> 282). This mirrors the
# OOXML diff: dimension grows from A1:R281 to A1:R282 and a brand-new record
# (fecha 44951) appears at row 246 while all subsequent rows keep their data
# but shift down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 246:281 down to 247:282, leaving a blank row 246 to fill in.
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row 246 with the new weekly record.
$ws.Cells.Item(246, 1).Value = 8
$ws.Cells.Item(246, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(246, 3).Value = "Coquimbo"
$ws.Cells.Item(246, 4).Value = 44951
$ws.Cells.Item(246, 5).Value = 4
$ws.Cells.Item(246, 6).Value = 100112037
$ws.Cells.Item(246, 7).Value = "Cebollín"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 1560
$ws.Cells.Item(246, 11).Value = 1200
$ws.Cells.Item(246, 12).Value = 1400
$ws.Cells.Item(246, 13).Value = 1300
$ws.Cells.Item(246, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(246, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(246, 16).Value = 217
$ws.Cells.Item(246, 17).Value = 6
$ws.Cells.Item(246, 18).Value = "Hortaliza"
